# Auto-generated Excel COM-interop edit script
# Applies the "Updated cryptos list" data refresh (Wed May 24 17:26:19 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.448.59'
$ws.Range("E2").Value = '  -3.08%  '
$ws.Range("D3").Value = '1.803.51'
$ws.Range("E3").Value = '  -2.82%  '
$ws.Range("E4").Value = '  +0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.006'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '307.99'
$ws.Range("E6").Value = '  -1.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4533'
$ws.Range("E7").Value = '  -1.66%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3643'
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07093'
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8725'
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07766'
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.29'
$ws.Range("E12").Value = '  -4.41%  '
$ws.Range("D13").Value = '1.816.74'
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.252'
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.325'
$ws.Range("E15").Value = '  -3.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.09'
$ws.Range("E16").Value = '  -5.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008558'
$ws.Range("E18").Value = '  -4.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").Value = '26.463.56'
$ws.Range("E20").Value = '  -3.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.21'
$ws.Range("E21").Value = '  -3.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.957'
$ws.Range("E22").Value = '  -3.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.38'
$ws.Range("E23").Value = '  -1.64%  '
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.33'
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.87'
$ws.Range("E26").Value = '  -2.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.987'
$ws.Range("E27").Value = '  -3.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.76'
$ws.Range("E28").Value = '  -2.77%  '
$ws.Range("E29").Value = '  -4.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08644'
$ws.Range("E30").Value = '  -2.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.030'
$ws.Range("E31").Value = '  -1.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7277'
$ws.Range("E32").Value = '  -6.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.435'
$ws.Range("E33").Value = '  -1.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.109'
$ws.Range("E34").Value = '  -5.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.004'
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.530'
$ws.Range("E36").Value = '  -8.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.072'
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01920'
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05066'
$ws.Range("E39").Value = '  -3.51%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.869'
$ws.Range("E40").Value = '  -2.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.941'
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4915'
$ws.Range("E42").Value = '  -4.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1565'
$ws.Range("E43").Value = '  -4.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.105'
$ws.Range("E44").Value = '  -3.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.007'
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4600'
$ws.Range("E46").Value = '  -4.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.41'
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.879'
$ws.Range("E48").Value = '  -4.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.579'
$ws.Range("E49").Value = '  -4.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.47'
$ws.Range("E51").Value = '  -3.40%  '
